# Regenerate orders with updated distance/size labels.
# The experiment condition codes change:
#   D64 -> D69, D80 -> D86, D51 -> D55, S30 -> S31
# These substitutions touch the "Condition", "Filename_Left", "Filename_Right",
# "Distance" and "Size" columns (and any other cell containing these tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cells = $ws.Cells

$cells.Replace("D64", "D69", 2, 1, $false, $false, $false, $false)
$cells.Replace("D80", "D86", 2, 1, $false, $false, $false, $false)
$cells.Replace("D51", "D55", 2, 1, $false, $false, $false, $false)
$cells.Replace("S30", "S31", 2, 1, $false, $false, $false, $false)
